# Adds a new column BB (dt_full_qoq_PUBCON_AVERAGE_10_9.xlsx update):
#   - BB1 gets a new quarter-end date header (45986)
#   - BB2:BB70 copy the prior column (BA) values forward one quarter
#   - BB71:BB83 carry new / revised forecast values
#   - A new row 83 is appended (date 46934) with its own BB83 value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 83 (extends column A down one more quarter) -----------------
$ws.Range("A83").Value = 46934
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122)   # xlPasteFormats - copy the date style

# --- New column BB ----------------------------------------------------
# Header date for the new column
$ws.Range("BB1").Value = 45986
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)   # xlPasteFormats - copy the header style

# BB2:BB70 simply repeat the BA column's values (no special cell style)
$ws.Range("BA2:BA70").Copy()
$ws.Range("BB2:BB70").PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = 0

# BB71:BB83 hold the newly revised forecast values
$ws.Range("BB71").Value = 0.1765865160815849
$ws.Range("BB72").Value = 0.2412052862208469
$ws.Range("BB73").Value = 0.768168485846715
$ws.Range("BB74").Value = 0.6732276193363788
$ws.Range("BB75").Value = 0.6732276193363788
$ws.Range("BB76").Value = 0.6732276193363788
$ws.Range("BB77").Value = 0.6732276193363788
$ws.Range("BB78").Value = 0.6732276193363788
$ws.Range("BB79").Value = 0.6732276193363788
$ws.Range("BB80").Value = 0.6732276193363788
$ws.Range("BB81").Value = 0.6732276193363788
$ws.Range("BB82").Value = 0.6732276193363788
$ws.Range("BB83").Value = 0.6732276193363788
